$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert 6 extra data rows (18:23) so the worker-debt table grows
#    from 3 rows (16-18) to 9 rows (16-24); this also pushes the
#    trailing signature block from rows 23-24 down to rows 29-30.
# ------------------------------------------------------------------
$ws.Rows("18:23").Insert()

# Copy the formatting (borders/fill/number format) of row 17 (a
# "middle" data row) onto the freshly inserted rows so they look like
# the rest of the table instead of picking up default styling.
$ws.Range("B17:J17").Copy() | Out-Null
$ws.Range("B18:J23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Header / summary block updates
# ------------------------------------------------------------------
$ws.Range("D2").Value = "ESTADO DE CUENTA"
$ws.Range("B7").Value = "RAZON SOCIAL:"
$ws.Range("B11").Value = "VALOR MORA"
$ws.Range("E11").Value = 362400

$ws.Range("B13").Value = "Cant. Trabajadores"
$ws.Range("C13").Value = 4
$ws.Range("E13").Value = "Cant. Periodos"
$ws.Range("F13").Value = 4

# Table header row
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"
$ws.Range("J15").Value = "Observaciones"

# ------------------------------------------------------------------
# 3. Worker debt rows (16-24)
# ------------------------------------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73155360"
$ws.Range("D16").Value = "DEWEY MAY BARRETO"
$ws.Range("E16").Value = "2304"
$ws.Range("F16").Value = 36000
$ws.Range("G16").Value = 1160000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73155360"
$ws.Range("D17").Value = "DEWEY MAY BARRETO"
$ws.Range("E17").Value = "2303"
$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 1160000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "73155360"
$ws.Range("D18").Value = "DEWEY MAY BARRETO"
$ws.Range("E18").Value = "2302"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1160000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73155360"
$ws.Range("D19").Value = "DEWEY MAY BARRETO"
$ws.Range("E19").Value = "2210"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1160000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "73111840"
$ws.Range("D20").Value = "WILMAR GUILLERMO ACUÑA BARRETO"
$ws.Range("E20").Value = "2302"
$ws.Range("F20").Value = 40000
$ws.Range("G20").Value = 0

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1047450067"
$ws.Range("D21").Value = "FABIO JOSE ZAPATA PACHECO"
$ws.Range("E21").Value = "2303"
$ws.Range("F21").Value = 40000
$ws.Range("G21").Value = 1500000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1047450067"
$ws.Range("D22").Value = "FABIO JOSE ZAPATA PACHECO"
$ws.Range("E22").Value = "2302"
$ws.Range("F22").Value = 40000
$ws.Range("G22").Value = 1500000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1047450067"
$ws.Range("D23").Value = "FABIO JOSE ZAPATA PACHECO"
$ws.Range("E23").Value = "2210"
$ws.Range("F23").Value = 40000
$ws.Range("G23").Value = 1500000

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1143385600"
$ws.Range("D24").Value = "ANDRES EDUARDO IBARRA LOZANO"
$ws.Range("E24").Value = "2302"
$ws.Range("F24").Value = 46400
$ws.Range("G24").Value = 1160000

# ------------------------------------------------------------------
# 4. Widen column D to fit the longest new worker name
# ------------------------------------------------------------------
$ws.Columns("D:D").ColumnWidth = 34.62

# ------------------------------------------------------------------
# 5. Footer block (previously rows 23-24, now pushed to 29-30 by the
#    row insert above) -- reassert the labels for safety.
# ------------------------------------------------------------------
$ws.Range("B29").Value = "___________________________________"
$ws.Range("H29").Value = "___________________________________"
$ws.Range("B30").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H30").Value = "FIRMA DEL REPRESENTANTE LEGAL"
